$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '40.106.32'
$ws.Cells.Item(2, 5).Value = '  -2.73%  '
$ws.Cells.Item(3, 4).Value = '2.341.08'
$ws.Cells.Item(3, 5).Value = '  -3.69%  '
$ws.Cells.Item(4, 5).Value = '  -0.15%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '309.86'
$ws.Cells.Item(5, 5).Value = '  -2.20%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '85.76'
$ws.Cells.Item(6, 5).Value = '  -3.38%  '
$ws.Cells.Item(7, 5).Value = '  -2.02%  '
$ws.Cells.Item(8, 5).Value = '  -0.04%  '
$ws.Cells.Item(9, 5).Value = '  -1.89%  '
$ws.Cells.Item(10, 5).Value = '  -2.36%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '30.10'
$ws.Cells.Item(11, 5).Value = '  -6.08%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.109'
$ws.Cells.Item(12, 5).Value = '  +0.94%  '
$ws.Cells.Item(13, 4).Value = '2.699.22'
$ws.Cells.Item(13, 5).Value = '  -3.90%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.43'
$ws.Cells.Item(14, 5).Value = '  -3.99%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '14.82'
$ws.Cells.Item(15, 5).Value = '  -4.59%  '
$ws.Cells.Item(16, 4).Value = '2.362.92'
$ws.Cells.Item(16, 5).Value = '  -3.65%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.758'
$ws.Cells.Item(17, 5).Value = '  -1.71%  '
$ws.Cells.Item(18, 4).Value = '40.084.14'
$ws.Cells.Item(18, 5).Value = '  -2.67%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0905'
$ws.Cells.Item(19, 5).Value = '  -1.94%  '
$ws.Cells.Item(20, 5).Value = '  -1.51%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '68.05'
$ws.Cells.Item(21, 5).Value = '  -4.92%  '
$ws.Cells.Item(22, 5).Value = '  -2.92%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '236.34'
$ws.Cells.Item(23, 5).Value = '  +0.64%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.55'
$ws.Cells.Item(24, 5).Value = '  -5.02%  '
$ws.Cells.Item(25, 5).Value = '  +0.23%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.83'
$ws.Cells.Item(26, 5).Value = '  -2.16%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '23.40'
$ws.Cells.Item(27, 5).Value = '  -2.24%  '
$ws.Cells.Item(28, 5).Value = '  -4.17%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.33'
$ws.Cells.Item(29, 5).Value = '  -2.16%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '35.10'
$ws.Cells.Item(30, 5).Value = '  +1.64%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '153.60'
$ws.Cells.Item(31, 5).Value = '  -2.25%  '
$ws.Cells.Item(32, 5).Value = '  -0.12%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.12'
$ws.Cells.Item(33, 5).Value = '  -2.57%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.44'
$ws.Cells.Item(34, 5).Value = '  -3.49%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0722'
$ws.Cells.Item(35, 5).Value = '  -2.68%  '
$ws.Cells.Item(36, 5).Value = '  -0.15%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.82'
$ws.Cells.Item(37, 5).Value = '  -2.51%  '
$ws.Cells.Item(38, 2).Value = 'Kaspa'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.0993'
$ws.Cells.Item(38, 5).Value = '  -0.20%  '
$ws.Cells.Item(39, 2).Value = 'Celestia'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '15.71'
$ws.Cells.Item(39, 5).Value = '  -4.82%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.72'
$ws.Cells.Item(40, 5).Value = '  -2.31%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.87'
$ws.Cells.Item(41, 5).Value = '  +0.31%  '
$ws.Cells.Item(42, 4).Value = '1.957.54'
$ws.Cells.Item(42, 5).Value = '  -1.16%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.22'
$ws.Cells.Item(43, 5).Value = '  -4.38%  '
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '17.63'
$ws.Cells.Item(44, 5).Value = '  -2.84%  '
$ws.Cells.Item(45, 2).Value = 'VeChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0264'
$ws.Cells.Item(45, 5).Value = '  -3.79%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '9.44'
$ws.Cells.Item(46, 5).Value = '  -0.68%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.71'
$ws.Cells.Item(47, 5).Value = '  -5.31%  '
$ws.Cells.Item(48, 4).Value = '2.559.35'
$ws.Cells.Item(48, 5).Value = '  -4.24%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '92.87'
$ws.Cells.Item(49, 5).Value = '  -2.41%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '70.85'
$ws.Cells.Item(50, 5).Value = '  -3.03%  '
$ws.Cells.Item(51, 2).Value = 'MultiversX'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '50.77'
$ws.Cells.Item(51, 5).Value = '  -2.19%  '

Write-Host "Applied 99 cell updates"
